# "Generate Report for Handoff"
#
# Localization status values flip from "In Translation" to "Ready for
# handoff" and the two "Latest *Datetime" timestamps that were mid-edit
# advance a few seconds/minutes. The three report sheets (Overview,
# zh-cn, de-de) all surface the shared "Status" text and/or the
# generation timestamp, so every occurrence on every sheet is updated.
# The Status/zh-cn/de-de columns also widen a bit to fit the new,
# longer "Ready for handoff" label.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$overview.Range("F2").Value = "Ready for handoff"   # de-de status column
$zhcn.Range("C2").Value     = "Ready for handoff"   # Status column
$dede.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Timestamps bumped forward by the new handoff generation run ------
$overview.Range("G2").Value = "2016-08-18 08:41:00" # Latest HO Xliff Generate Date
$dede.Range("H2").Value     = "2016-08-18 08:41:00" # Latest Handoff Datetime (de-de)
$zhcn.Range("H2").Value     = "2016-08-18 08:40:54" # Latest Handoff Datetime (zh-cn)

# --- Widen the Status / zh-cn / de-de columns to fit the longer label -
# Excel's ColumnWidth (character units) snaps to the nearest on-screen
# pixel, same as a live resize in the UI; 16 1/3 characters is the
# closest setting to the target ~17.22-character stored width.
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333  # E: zh-cn
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333  # F: de-de
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333333333333  # C: Status
$dede.Columns.Item(3).ColumnWidth     = 16.3333333333333  # C: Status
